# Apply the edit: clear the "200_120" header cell and correct the
# fusion ggnet_pretained/vgg 100_120 score.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 used to hold the "200_120" column header - that column/label is being
# removed from the shared-string table, so the header cell becomes blank
# (keeping its existing style).
$ws.Range("D1").ClearContents()

# D6 (fusion ggnet_pretained et vgg / 100_120) value corrected 0.85 -> 0.87
$ws.Range("D6").Value = 0.87

# Move the active selection to D9 (matches the saved selection in the file).
$ws.Range("D9").Select()
